$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Neodymium")
$ws1.Range("C1").Value = 2030
$ws1.Range("C2").Value = [double]"2.438729849902814E-05"
$ws1.Range("D2").Value = [double]"0.1236014393403301"
$ws1.Range("E2").Value = [double]"0.584289849462961"
$ws1.Range("B3").Value = [double]"2.183980475909259E-12"
$ws1.Range("C3").Value = [double]"0.001187950236302112"
$ws1.Range("D3").Value = [double]"0.1155224455664189"
$ws1.Range("E3").Value = [double]"0.4948953125759609"
$ws1.Range("B4").Value = [double]"3.409259119931335E-14"
$ws1.Range("C4").Value = [double]"0.001073104153047114"
$ws1.Range("D4").Value = [double]"0.08248514375558298"
$ws1.Range("E4").Value = [double]"0.4143225657740967"
$ws1.Range("C5").Value = [double]"2.380640363208162E-08"
$ws1.Range("D5").Value = [double]"0.004220506989166603"
$ws1.Range("E5").Value = [double]"0.03292316497354641"

$ws2 = $wb.Worksheets.Item("Dysprosium")
$ws2.Range("C1").Value = 2030
$ws2.Range("C2").Value = [double]"2.763358668632657E-05"
$ws2.Range("D2").Value = [double]"0.1067081964413131"
$ws2.Range("E2").Value = [double]"0.6620669446318078"
$ws2.Range("C3").Value = [double]"0.001346082914235243"
$ws2.Range("D3").Value = [double]"0.0997334001988444"
$ws2.Range("E3").Value = [double]"0.5607727531308762"
$ws2.Range("C4").Value = [double]"0.001215949222004488"
$ws2.Range("D4").Value = [double]"0.07121147593698561"
$ws2.Range("E4").Value = [double]"0.4694746545164075"
$ws2.Range("C5").Value = [double]"2.697536664354322E-08"
$ws2.Range("D5").Value = [double]"0.003643668644034803"
$ws2.Range("E5").Value = [double]"0.03730569555791484"

$ws3 = $wb.Worksheets.Item("Copper")
$ws3.Range("C1").Value = 2030
$ws3.Range("B2").Value = [double]"3.278497091721097E-06"
$ws3.Range("C2").Value = [double]"0.003050246220774824"
$ws3.Range("D2").Value = [double]"0.8588525024148915"
$ws3.Range("E2").Value = [double]"0.9697180417558001"
$ws3.Range("B3").Value = [double]"2.229370101113288E-05"
$ws3.Range("C3").Value = [double]"0.01103680953317707"
$ws3.Range("D3").Value = [double]"0.6168778258139214"
$ws3.Range("E3").Value = [double]"0.6805195557530022"
$ws3.Range("B4").Value = [double]"6.612099022439717E-05"
$ws3.Range("C4").Value = [double]"0.002936712507453067"
$ws3.Range("D4").Value = [double]"0.4416445656454649"
$ws3.Range("E4").Value = [double]"0.6055420661490584"
$ws3.Range("B5").Value = [double]"2.076994439830034E-05"
$ws3.Range("C5").Value = [double]"0.006488512101692222"
$ws3.Range("D5").Value = [double]"0.8181956015862759"
$ws3.Range("E5").Value = [double]"0.7122850722825392"

$ws4 = $wb.Worksheets.Item("Raw silicon")
$ws4.Range("C1").Value = 2030
$ws4.Range("B2").Value = [double]"4.96652837099915E-07"
$ws4.Range("C2").Value = [double]"0.0005194103187024707"
$ws4.Range("D2").Value = [double]"0.4783497257479737"
$ws4.Range("E2").Value = [double]"1.251484198240786"
$ws4.Range("B3").Value = [double]"5.30035999530297E-07"
$ws4.Range("C3").Value = [double]"0.001746085699123397"
$ws4.Range("D3").Value = [double]"0.2315523203444246"
$ws4.Range("E3").Value = [double]"0.5836018294026056"
$ws4.Range("B4").Value = [double]"3.397047964529607E-06"
$ws4.Range("C4").Value = [double]"0.0004869161330086471"
$ws4.Range("D4").Value = [double]"0.2164427902680356"
$ws4.Range("E4").Value = [double]"0.6322978441499242"
$ws4.Range("B5").Value = [double]"1.823860200208514E-06"
$ws4.Range("C5").Value = [double]"0.0006186310519116159"
$ws4.Range("D5").Value = [double]"0.4280362603778582"
$ws4.Range("E5").Value = [double]"0.8793987462358255"
